$wb = $excel.ActiveWorkbook

# --- Sheet: steel_prim ---
$ws = $wb.Worksheets.Item("steel_prim")
$ws.Range("D34").Value = "(2018.0, 0.0)"
$ws.Range("G34").Value = 0.0
$ws.Range("H34").Value = 9.96

# --- Sheet: steel_sec ---
$ws = $wb.Worksheets.Item("steel_sec")
$ws.Range("F4").Value = 0.0
$ws.Range("H4").Value = -0.91
$ws.Range("I4").Value = 0.0
$ws.Range("J4").Value = -0.0
$ws.Range("F6").Value = 0.01
$ws.Range("G6").Value = -0.0
$ws.Range("H6").Value = 0.98
$ws.Range("H8").Value = -8.97
$ws.Range("F10").Value = 0.01
$ws.Range("H10").Value = 3.87
$ws.Range("I10").Value = -0.0
$ws.Range("J10").Value = 0.0
$ws.Range("H11").Value = 1.9
$ws.Range("I11").Value = -0.0
$ws.Range("J11").Value = 0.0
$ws.Range("F12").Value = 0.0
$ws.Range("G12").Value = -0.0
$ws.Range("H12").Value = 0.59
$ws.Range("I12").Value = -0.0
$ws.Range("J12").Value = 0.0
$ws.Range("F15").Value = -0.01
$ws.Range("H15").Value = 18.67
$ws.Range("I15").Value = -0.02
$ws.Range("F18").Value = -0.02
$ws.Range("H18").Value = 3.92
$ws.Range("F19").Value = -0.02
$ws.Range("H19").Value = -1.73
$ws.Range("H21").Value = 191.17
$ws.Range("I21").Value = -0.19
$ws.Range("J21").Value = 0.0
$ws.Range("H22").Value = -1.51
$ws.Range("I22").Value = 0.0
$ws.Range("J22").Value = -0.0
$ws.Range("H25").Value = -1.79
$ws.Range("H27").Value = -8.51
$ws.Range("I27").Value = 0.01
$ws.Range("C29").Value = "(2018.0, 0.0001281860992295341)"
$ws.Range("H31").Value = 1.95
$ws.Range("I31").Value = -0.0
$ws.Range("J31").Value = 0.0

# --- Sheet: alu_prim ---
$ws = $wb.Worksheets.Item("alu_prim")
$ws.Range("F6").Value = -0.0
$ws.Range("G6").Value = 0.0
$ws.Range("H6").Value = 0.98
$ws.Range("H8").Value = -8.97
$ws.Range("F10").Value = -0.0
$ws.Range("G10").Value = 0.0
$ws.Range("H10").Value = 3.87
$ws.Range("I10").Value = -0.0
$ws.Range("J10").Value = 0.0
$ws.Range("F21").Value = -0.0
$ws.Range("G21").Value = 0.0
$ws.Range("H21").Value = 191.17
$ws.Range("I21").Value = -0.19
$ws.Range("J21").Value = 0.0
$ws.Range("H22").Value = -1.51
$ws.Range("I22").Value = 0.0
$ws.Range("J22").Value = -0.0
$ws.Range("H25").Value = -1.79

# --- Sheet: chlorine ---
$ws = $wb.Worksheets.Item("chlorine")
$ws.Range("H4").Value = -0.91
$ws.Range("I4").Value = 0.0
$ws.Range("J4").Value = -0.0
$ws.Range("F6").Value = 0.0
$ws.Range("G6").Value = -0.0
$ws.Range("H6").Value = 0.98
$ws.Range("H7").Value = -112.99
$ws.Range("I7").Value = 0.11
$ws.Range("H8").Value = -8.97
$ws.Range("H10").Value = 3.87
$ws.Range("I10").Value = -0.0
$ws.Range("J10").Value = 0.0
$ws.Range("H12").Value = 0.59
$ws.Range("I12").Value = -0.0
$ws.Range("J12").Value = 0.0
$ws.Range("F15").Value = -0.01
$ws.Range("H15").Value = 18.67
$ws.Range("I15").Value = -0.02
$ws.Range("F18").Value = -0.0
$ws.Range("G18").Value = 0.0
$ws.Range("H18").Value = 3.92
$ws.Range("H19").Value = -1.73
$ws.Range("H21").Value = 191.17
$ws.Range("I21").Value = -0.19
$ws.Range("J21").Value = 0.0
$ws.Range("H22").Value = -1.51
$ws.Range("I22").Value = 0.0
$ws.Range("J22").Value = -0.0
$ws.Range("H25").Value = -1.79
$ws.Range("H27").Value = -8.51
$ws.Range("I27").Value = 0.01

# --- Sheet: paper ---
$ws = $wb.Worksheets.Item("paper")
$ws.Range("H4").Value = -0.91
$ws.Range("I4").Value = 0.0
$ws.Range("J4").Value = -0.0
$ws.Range("H5").Value = 74.25
$ws.Range("I5").Value = -0.07
$ws.Range("J5").Value = 0.0
$ws.Range("F6").Value = 0.0
$ws.Range("G6").Value = -0.0
$ws.Range("H6").Value = 0.98
$ws.Range("H7").Value = -112.99
$ws.Range("I7").Value = 0.11
$ws.Range("H8").Value = -8.97
$ws.Range("H10").Value = 3.87
$ws.Range("I10").Value = -0.0
$ws.Range("J10").Value = 0.0
$ws.Range("F11").Value = 0.0
$ws.Range("G11").Value = -0.0
$ws.Range("H11").Value = 1.9
$ws.Range("I11").Value = -0.0
$ws.Range("J11").Value = 0.0
$ws.Range("F12").Value = 0.0
$ws.Range("G12").Value = -0.0
$ws.Range("H12").Value = 0.59
$ws.Range("I12").Value = -0.0
$ws.Range("J12").Value = 0.0
$ws.Range("F15").Value = -0.04
$ws.Range("H15").Value = 18.67
$ws.Range("I15").Value = -0.02
$ws.Range("H18").Value = 3.92
$ws.Range("H19").Value = -1.73
$ws.Range("F21").Value = -0.0
$ws.Range("H21").Value = 191.17
$ws.Range("I21").Value = -0.19
$ws.Range("J21").Value = 0.0
$ws.Range("F22").Value = -0.0
$ws.Range("H22").Value = -1.51
$ws.Range("I22").Value = 0.0
$ws.Range("J22").Value = -0.0
$ws.Range("H25").Value = -1.79
$ws.Range("H27").Value = -8.51
$ws.Range("I27").Value = 0.01
$ws.Range("C29").Value = "(2018.0, 1.1565663088379015e-05)"
$ws.Range("F31").Value = -0.0
$ws.Range("H31").Value = 1.95
$ws.Range("I31").Value = -0.0
$ws.Range("J31").Value = 0.0
$ws.Range("H34").Value = 9.96

# --- Sheet: cement ---
$ws = $wb.Worksheets.Item("cement")
$ws.Range("F4").Value = 0.01
$ws.Range("H4").Value = -0.91
$ws.Range("I4").Value = 0.0
$ws.Range("J4").Value = -0.0
$ws.Range("H5").Value = 74.25
$ws.Range("I5").Value = -0.07
$ws.Range("J5").Value = 0.0
$ws.Range("H6").Value = 0.98
$ws.Range("F7").Value = -0.11
$ws.Range("H7").Value = -112.99
$ws.Range("I7").Value = 0.11
$ws.Range("F8").Value = 0.06
$ws.Range("H8").Value = -8.97
$ws.Range("H10").Value = 3.87
$ws.Range("I10").Value = -0.0
$ws.Range("J10").Value = 0.0
$ws.Range("F11").Value = 0.04
$ws.Range("G11").Value = -0.0
$ws.Range("H11").Value = 1.9
$ws.Range("I11").Value = -0.0
$ws.Range("J11").Value = 0.0
$ws.Range("F12").Value = 0.06
$ws.Range("H12").Value = 0.59
$ws.Range("I12").Value = -0.0
$ws.Range("J12").Value = 0.0
$ws.Range("F15").Value = -0.01
$ws.Range("G15").Value = 0.0
$ws.Range("H15").Value = 18.67
$ws.Range("I15").Value = -0.02
$ws.Range("H18").Value = 3.92
$ws.Range("F19").Value = 0.04
$ws.Range("H19").Value = -1.73
$ws.Range("F21").Value = 0.02
$ws.Range("H21").Value = 191.17
$ws.Range("I21").Value = -0.19
$ws.Range("J21").Value = 0.0
$ws.Range("H22").Value = -1.51
$ws.Range("I22").Value = 0.0
$ws.Range("J22").Value = -0.0
$ws.Range("F25").Value = -0.01
$ws.Range("G25").Value = 0.0
$ws.Range("H25").Value = -1.79
$ws.Range("F27").Value = 0.0
$ws.Range("H27").Value = -8.51
$ws.Range("I27").Value = 0.01
$ws.Range("C29").Value = "(2018.0, 0.0004515950085971429)"
$ws.Range("F30").Value = 0.02
$ws.Range("G30").Value = -0.0
$ws.Range("H30").Value = -17.15
$ws.Range("I30").Value = 0.02
$ws.Range("F31").Value = -0.01
$ws.Range("G31").Value = 0.0
$ws.Range("H31").Value = 1.95
$ws.Range("I31").Value = -0.0
$ws.Range("J31").Value = 0.0
$ws.Range("F34").Value = -0.02
$ws.Range("H34").Value = 9.96
